# Re-applies the scheduled runner's recalculated Leve profit figures
# (currentAveragePrice* / LevePrice* / LeveProfit* columns, H:N) across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 419.8
$ws.Range("I33").Value = 299.66666
$ws.Range("K33").Value = 299.66666
$ws.Range("M33").Value = -70.66665999999998
$ws.Range("H40").Value = 2229.35
$ws.Range("I40").Value = 2084.3333
$ws.Range("K40").Value = 2084.3333
$ws.Range("M40").Value = -1909.3333
$ws.Range("H69").Value = 32714.928
$ws.Range("J69").Value = 19166.416
$ws.Range("L69").Value = 57499.24800000001
$ws.Range("N69").Value = -59247.24800000001
$ws.Range("H72").Value = 32714.928
$ws.Range("J72").Value = 19166.416
$ws.Range("L72").Value = 172497.744
$ws.Range("N72").Value = -181233.744
$ws.Range("H100").Value = 1177
$ws.Range("I100").Value = 221.25
$ws.Range("J100").Value = 5000
$ws.Range("K100").Value = 221.25
$ws.Range("L100").Value = 5000
$ws.Range("M100").Value = 319.75
$ws.Range("N100").Value = -6082
$ws.Range("H137").Value = 2428.875
$ws.Range("I137").Value = 2173.7144
$ws.Range("K137").Value = 6521.1432
$ws.Range("M137").Value = -3971.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1062
$ws.Range("I45").Value = 1062
$ws.Range("K45").Value = 1062
$ws.Range("M45").Value = -685
$ws.Range("H61").Value = 4995.353
$ws.Range("I61").Value = 1616.5
$ws.Range("K61").Value = 1616.5
$ws.Range("M61").Value = -1404.5
$ws.Range("H63").Value = 3750
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 4500
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 4500
$ws.Range("M63").Value = -2314
$ws.Range("N63").Value = -5872
$ws.Range("H66").Value = 3750
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 4500
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 22500
$ws.Range("M66").Value = -11568
$ws.Range("N66").Value = -29364
$ws.Range("H74").Value = 2559.48
$ws.Range("I74").Value = 1908.625
$ws.Range("J74").Value = 3716.5557
$ws.Range("K74").Value = 1908.625
$ws.Range("L74").Value = 3716.5557
$ws.Range("M74").Value = -1034.625
$ws.Range("N74").Value = -5464.5557
$ws.Range("H77").Value = 2559.48
$ws.Range("I77").Value = 1908.625
$ws.Range("J77").Value = 3716.5557
$ws.Range("K77").Value = 9543.125
$ws.Range("L77").Value = 18582.7785
$ws.Range("M77").Value = -5175.125
$ws.Range("N77").Value = -27318.7785
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H122").Value = 1959.6666
$ws.Range("I122").Value = 1449.75
$ws.Range("J122").Value = 2542.4285
$ws.Range("K122").Value = 4349.25
$ws.Range("L122").Value = 7627.2855
$ws.Range("M122").Value = -1899.25
$ws.Range("N122").Value = -12527.2855
$ws.Range("H136").Value = 4995.353
$ws.Range("I136").Value = 1616.5
$ws.Range("K136").Value = 4849.5
$ws.Range("M136").Value = -2299.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1294.5
$ws.Range("I20").Value = 1294.5
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1294.5
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1047.5
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 771.4375
$ws.Range("I22").Value = 588.5833
$ws.Range("K22").Value = 588.5833
$ws.Range("M22").Value = -415.5833
$ws.Range("H94").Value = 6608.5454
$ws.Range("I94").Value = 7159.5
$ws.Range("K94").Value = 7159.5
$ws.Range("M94").Value = -6708.5
$ws.Range("H99").Value = 2276.6667
$ws.Range("I99").Value = 1450.1111
$ws.Range("J99").Value = 3516.5
$ws.Range("K99").Value = 1450.1111
$ws.Range("L99").Value = 3516.5
$ws.Range("M99").Value = 47.88889999999992
$ws.Range("N99").Value = -6512.5
$ws.Range("H107").Value = 4050
$ws.Range("I107").Value = 3733.3333
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 3733.3333
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = -1813.3333
$ws.Range("N107").Value = -8840
$ws.Range("H134").Value = 1311.3158
$ws.Range("I134").Value = 1311.3158
$ws.Range("K134").Value = 3933.9474
$ws.Range("M134").Value = -1398.9474

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2761.111
$ws.Range("I31").Value = 2226.5386
$ws.Range("K31").Value = 2226.5386
$ws.Range("M31").Value = -1931.5386
$ws.Range("H34").Value = 2761.111
$ws.Range("I34").Value = 2226.5386
$ws.Range("K34").Value = 2226.5386
$ws.Range("M34").Value = -2024.5386
$ws.Range("H86").Value = 19950
$ws.Range("I86").Value = 19950
$ws.Range("K86").Value = 19950
$ws.Range("M86").Value = -18827
$ws.Range("H89").Value = 19950
$ws.Range("I89").Value = 19950
$ws.Range("K89").Value = 99750
$ws.Range("M89").Value = -94134
$ws.Range("H122").Value = 1890.1666
$ws.Range("I122").Value = 1538.3334
$ws.Range("K122").Value = 4615.0002
$ws.Range("M122").Value = -2165.0002
$ws.Range("H132").Value = 2883.1
$ws.Range("I132").Value = 1558.25
$ws.Range("K132").Value = 4674.75
$ws.Range("M132").Value = -2144.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1133.3334
$ws.Range("I18").Value = 1133.3334
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 3400.0002
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -3231.0002
$ws.Range("N18").ClearContents()
$ws.Range("H44").Value = 2000
$ws.Range("I44").Value = 2000
$ws.Range("K44").Value = 6000
$ws.Range("M44").Value = -5602
$ws.Range("H58").Value = 2196.8
$ws.Range("J58").Value = 2246.25
$ws.Range("L58").Value = 6738.75
$ws.Range("N58").Value = -6994.75
$ws.Range("H80").Value = 2294.6667
$ws.Range("J80").Value = 2294.6667
$ws.Range("L80").Value = 6884.000100000001
$ws.Range("N80").Value = -8756.000100000001
$ws.Range("H83").Value = 2294.6667
$ws.Range("J83").Value = 2294.6667
$ws.Range("L83").Value = 20652.0003
$ws.Range("N83").Value = -30012.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 29507
$ws.Range("J47").Value = 29507
$ws.Range("L47").Value = 29507
$ws.Range("N47").Value = -30643
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()
$ws.Range("H102").Value = 1536.9231
$ws.Range("I102").Value = 1028
$ws.Range("J102").Value = 3233.3333
$ws.Range("K102").Value = 1028
$ws.Range("L102").Value = 3233.3333
$ws.Range("M102").Value = 594
$ws.Range("N102").Value = -6477.3333
$ws.Range("H113").Value = 1175.8889
$ws.Range("I113").Value = 1098
$ws.Range("J113").Value = 1448.5
$ws.Range("K113").Value = 1098
$ws.Range("L113").Value = 1448.5
$ws.Range("M113").Value = 1072
$ws.Range("N113").Value = -5788.5
$ws.Range("H132").Value = 2971.6667
$ws.Range("I132").Value = 2971.6667
$ws.Range("K132").Value = 8915.000100000001
$ws.Range("M132").Value = -6385.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 22423.25
$ws.Range("J46").Value = 1483.421
$ws.Range("L46").Value = 1483.421
$ws.Range("N46").Value = -1859.421
$ws.Range("H82").Value = 2001.3846
$ws.Range("I82").Value = 1513.3334
$ws.Range("K82").Value = 1513.3334
$ws.Range("M82").Value = -1152.3334
$ws.Range("H85").Value = 2001.3846
$ws.Range("I85").Value = 1513.3334
$ws.Range("K85").Value = 1513.3334
$ws.Range("M85").Value = -265.3334
$ws.Range("H100").Value = 4833.3335
$ws.Range("I100").Value = 3000
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459
$ws.Range("H122").Value = 6742.613
$ws.Range("I122").Value = 7106.8945
$ws.Range("K122").Value = 21320.6835
$ws.Range("M122").Value = -18870.6835
$ws.Range("H132").Value = 3410.6
$ws.Range("I132").Value = 1274
$ws.Range("J132").Value = 4835
$ws.Range("K132").Value = 3822
$ws.Range("L132").Value = 14505
$ws.Range("M132").Value = -1292
$ws.Range("N132").Value = -19565

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1458.125
$ws.Range("I107").Value = 882.5
$ws.Range("J107").Value = 1650
$ws.Range("K107").Value = 2647.5
$ws.Range("L107").Value = 4950
$ws.Range("M107").Value = -727.5
$ws.Range("N107").Value = -8790
$ws.Range("H132").Value = 1872.7142
$ws.Range("I132").Value = 1778.25
$ws.Range("J132").Value = 1998.6666
$ws.Range("K132").Value = 5334.75
$ws.Range("L132").Value = 5995.9998
$ws.Range("M132").Value = -2804.75
$ws.Range("N132").Value = -11055.9998
$ws.Range("H136").Value = 3597.8823
$ws.Range("I136").Value = 2077.7334
$ws.Range("K136").Value = 6233.2002
$ws.Range("M136").Value = -3683.2002
$ws.Range("H137").Value = 66666
$ws.Range("J137").Value = 66666
$ws.Range("L137").Value = 66666
$ws.Range("N137").Value = -76866

Write-Host "Applied 247 cell edits (239 updated, 8 cleared) across 8 sheets"
